# Adding PC SunEnergy to the Portfolio Forecast
# Shift all timestamps in column A by +2 days (new data window) and
# update the "Notified Production (MW)" values in column B to reflect
# the merged/aggregated portfolio totals that now include PC SunEnergy.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 97

# Shift every timestamp in column A down by 2 days (keeps same time-of-day
# fraction, just moves the date from 2025-05-20/21 to 2025-05-22/23).
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $cell.Value2 + 2
}

# New "Notified Production (MW)" values for the rows whose totals changed
# (row number -> new value).
$newB = @{
    19 = 1;    20 = 1;    21 = 2;    22 = 47;   23 = 54;   24 = 64;   25 = 79;
    26 = 310;  27 = 333;  28 = 359;  29 = 389;  30 = 762;  31 = 794;  32 = 802;
    33 = 805;  34 = 1112; 35 = 1150; 36 = 1166; 37 = 1215; 38 = 1436; 39 = 1466;
    40 = 1502; 41 = 1528; 42 = 1662; 43 = 1673; 44 = 1683; 45 = 1691; 46 = 1749;
    47 = 1757; 48 = 1761; 49 = 1762; 50 = 1780; 51 = 1774; 52 = 1766; 53 = 1754;
    54 = 1631; 55 = 1604; 56 = 1584; 57 = 1564; 58 = 1418; 59 = 1387; 60 = 1366;
    61 = 1339; 62 = 1144; 63 = 1127; 64 = 1103; 65 = 1084; 66 = 877;  67 = 857;
    68 = 833;  69 = 811;  70 = 531;  71 = 510;  72 = 487;  73 = 460;  74 = 199;
    75 = 179;  76 = 157;  77 = 141;  78 = 34;   79 = 27;   80 = 22;   81 = 21;
    83 = 36;   84 = 36;   85 = 36;   86 = 1;    87 = 1;    88 = 1;    89 = 1
}

foreach ($row in $newB.Keys) {
    $ws.Cells.Item([int]$row, 2).Value = $newB[$row]
}
